$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.514.25'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.376.21'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.93%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.67'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.545'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.386.05'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.08%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.84%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.330'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.94%  '

$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.84'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.798.79'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.442.04'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.63'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.99%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.404.31'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '309.50'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.28'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.35%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.23'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.14%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.10%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.58'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.36%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.46%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.50%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.995'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.67'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.62%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.99%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.34%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.47'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.58%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '127.37'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.77'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.565'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '240.06'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0483'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.05'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.67%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.20%  '
